$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): split each "<group>" header into "<group> mean" / "<group> std" ---
$ws.Range("B1").Value = "Algorithm"
$ws.Range("C1").Value = "State Based mean"
$ws.Range("D1").Value = "State Based std"
$ws.Range("E1").Value = "Non State mean"
$ws.Range("F1").Value = "Non State std"
$ws.Range("G1").Value = "One Sided mean"
$ws.Range("H1").Value = "One Sided std"

# Give the two newly-added header cells (G1, H1 reuse existing column style already;
# F1/G1/H1 are brand new columns) the same look (bold / bordered / centered) as the
# rest of the header row.
$ws.Range("B1").Copy() | Out-Null
foreach ($addr in @("F1", "G1", "H1")) {
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false

# --- Data rows: algorithm names (col B) and the mean/std values (cols C-H) ---
$algorithms = @("LR", "LDA", "KNN", "DTREE", "RTREE", "XTREE", "SVM")

$values = @(
    @(0.8829013906447534, 0.03274199378815969, 0.8711213517665131, 0.04629597460331539, 0.8677814845704754, 0.02473225417214233),
    @(0.8820243362831859, 0.03463622964191124, 0.8855606758832565, 0.05133304089556665, 0.8650792326939116, 0.0194877306875019),
    @(0.899778761061947,  0.03348585117764191, 0.8598310291858677, 0.03837979361932545, 0.8687072560467055, 0.03088087950425139),
    @(0.7720290771175727, 0.03634337026368784, 0.6763696876600103, 0.05243368698360654, 0.7365304420350292, 0.03467672869555064),
    @(0.8828934892541087, 0.02080682216664726, 0.8936763952892985, 0.03770871062676921, 0.893369474562135,  0.02255618088378395),
    @(0.8882269279393172, 0.03356348919724214, 0.8888376856118791, 0.03857748133936107, 0.877839866555463,  0.01881406225487286),
    @(0.9015170670037925, 0.01872160108006543, 0.8904505888376855, 0.04012577235616578, 0.8714762301918263, 0.02416475367188197)
)

for ($i = 0; $i -lt 7; $i++) {
    $row = $i + 2

    $ws.Cells.Item($row, 2).Value = $algorithms[$i]

    $rowVals = $values[$i]
    for ($j = 0; $j -lt 6; $j++) {
        $ws.Cells.Item($row, $j + 3).Value = $rowVals[$j]
    }
}

# The table used to have one more row (index 7, algorithm "NB") which no longer
# exists - clear it out entirely so the used range shrinks back to A1:H8.
$ws.Range("A9:H9").Clear() | Out-Null
